$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry updates one cell on the "cryptos" sheet to the latest scraped
# value. Cells whose text looks like a plain decimal number (e.g. "263.09")
# are explicitly formatted as Text ("@") before the assignment so Excel
# keeps them as literal strings (matching the source data, which stores
# prices as text) instead of silently re-interpreting them as numbers.
$updates = @(
    @{Cell='D2'; Text='26.493.54'; Numeric=$False},
    @{Cell='E2'; Text='  +0.02%  '; Numeric=$False},
    @{Cell='D3'; Text='1.844.61'; Numeric=$False},
    @{Cell='E3'; Text='  -0.37%  '; Numeric=$False},
    @{Cell='E4'; Text='  -0.01%  '; Numeric=$False},
    @{Cell='D5'; Text='263.09'; Numeric=$True},
    @{Cell='E5'; Text='  +0.63%  '; Numeric=$False},
    @{Cell='E6'; Text='  +0.00%  '; Numeric=$False},
    @{Cell='D7'; Text='0.5214'; Numeric=$True},
    @{Cell='E7'; Text='  +1.32%  '; Numeric=$False},
    @{Cell='D8'; Text='0.3229'; Numeric=$True},
    @{Cell='E8'; Text='  -1.36%  '; Numeric=$False},
    @{Cell='D9'; Text='0.06787'; Numeric=$True},
    @{Cell='E9'; Text='  +0.12%  '; Numeric=$False},
    @{Cell='D10'; Text='18.67'; Numeric=$True},
    @{Cell='E10'; Text='  -2.19%  '; Numeric=$False},
    @{Cell='D11'; Text='0.7748'; Numeric=$True},
    @{Cell='E11'; Text='  +0.03%  '; Numeric=$False},
    @{Cell='E12'; Text='  +0.94%  '; Numeric=$False},
    @{Cell='D13'; Text='1.859.14'; Numeric=$False},
    @{Cell='E13'; Text='  +0.69%  '; Numeric=$False},
    @{Cell='D14'; Text='88.27'; Numeric=$True},
    @{Cell='E14'; Text='  -0.49%  '; Numeric=$False},
    @{Cell='D15'; Text='5.011'; Numeric=$True},
    @{Cell='E15'; Text='  -0.61%  '; Numeric=$False},
    @{Cell='D16'; Text='1.000'; Numeric=$True},
    @{Cell='E16'; Text='  -0.02%  '; Numeric=$False},
    @{Cell='D17'; Text='13.94'; Numeric=$True},
    @{Cell='E17'; Text='  -1.48%  '; Numeric=$False},
    @{Cell='B18'; Text='ShibaInu'; Numeric=$False},
    @{Cell='C18'; Text='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; Numeric=$False},
    @{Cell='D18'; Text='0.000007965'; Numeric=$True},
    @{Cell='E18'; Text='  +0.65%  '; Numeric=$False},
    @{Cell='B19'; Text='Dai'; Numeric=$False},
    @{Cell='C19'; Text='https://coinranking.com/coin/MoTuySvg7+dai-dai'; Numeric=$False},
    @{Cell='D19'; Text='1.000'; Numeric=$True},
    @{Cell='E19'; Text='  +0.06%  '; Numeric=$False},
    @{Cell='D20'; Text='26.536.40'; Numeric=$False},
    @{Cell='E20'; Text='  +0.05%  '; Numeric=$False},
    @{Cell='D21'; Text='2.088.45'; Numeric=$False},
    @{Cell='E21'; Text='  +0.52%  '; Numeric=$False},
    @{Cell='D22'; Text='4.612'; Numeric=$True},
    @{Cell='E22'; Text='  +1.22%  '; Numeric=$False},
    @{Cell='D23'; Text='9.443'; Numeric=$True},
    @{Cell='E23'; Text='  -1.15%  '; Numeric=$False},
    @{Cell='D24'; Text='5.982'; Numeric=$True},
    @{Cell='E24'; Text='  +0.40%  '; Numeric=$False},
    @{Cell='D25'; Text='142.84'; Numeric=$True},
    @{Cell='E25'; Text='  -1.25%  '; Numeric=$False},
    @{Cell='D26'; Text='2.162'; Numeric=$True},
    @{Cell='E26'; Text='  -8.66%  '; Numeric=$False},
    @{Cell='D27'; Text='1.681'; Numeric=$True},
    @{Cell='E27'; Text='  +1.56%  '; Numeric=$False},
    @{Cell='D28'; Text='16.99'; Numeric=$True},
    @{Cell='E28'; Text='  +0.18%  '; Numeric=$False},
    @{Cell='D29'; Text='111.85'; Numeric=$True},
    @{Cell='E29'; Text='  +0.49%  '; Numeric=$False},
    @{Cell='D30'; Text='4.158'; Numeric=$True},
    @{Cell='E30'; Text='  -1.39%  '; Numeric=$False},
    @{Cell='D31'; Text='0.08730'; Numeric=$True},
    @{Cell='D32'; Text='4.103'; Numeric=$True},
    @{Cell='E32'; Text='  -1.58%  '; Numeric=$False},
    @{Cell='D33'; Text='0.04816'; Numeric=$True},
    @{Cell='E33'; Text='  -0.88%  '; Numeric=$False},
    @{Cell='D34'; Text='0.7185'; Numeric=$True},
    @{Cell='E34'; Text='  +3.72%  '; Numeric=$False},
    @{Cell='E35'; Text='  -1.05%  '; Numeric=$False},
    @{Cell='D36'; Text='2.861'; Numeric=$True},
    @{Cell='E36'; Text='  +0.65%  '; Numeric=$False},
    @{Cell='E37'; Text='  -1.06%  '; Numeric=$False},
    @{Cell='D38'; Text='0.01792'; Numeric=$True},
    @{Cell='E38'; Text='  -0.84%  '; Numeric=$False},
    @{Cell='D39'; Text='2.202'; Numeric=$True},
    @{Cell='E39'; Text='  -0.92%  '; Numeric=$False},
    @{Cell='E40'; Text='  -1.92%  '; Numeric=$False},
    @{Cell='D41'; Text='111.19'; Numeric=$True},
    @{Cell='E41'; Text='  -2.16%  '; Numeric=$False},
    @{Cell='D42'; Text='0.8876'; Numeric=$True},
    @{Cell='E42'; Text='  -0.99%  '; Numeric=$False},
    @{Cell='D43'; Text='6.037'; Numeric=$True},
    @{Cell='E43'; Text='  -1.85%  '; Numeric=$False},
    @{Cell='D44'; Text='0.9998'; Numeric=$True},
    @{Cell='E44'; Text='  +0.01%  '; Numeric=$False},
    @{Cell='D45'; Text='7.607'; Numeric=$True},
    @{Cell='E45'; Text='  -2.57%  '; Numeric=$False},
    @{Cell='E46'; Text='  -2.15%  '; Numeric=$False},
    @{Cell='D47'; Text='0.05891'; Numeric=$True},
    @{Cell='E47'; Text='  -0.15%  '; Numeric=$False},
    @{Cell='D48'; Text='9.049'; Numeric=$True},
    @{Cell='E48'; Text='  -0.77%  '; Numeric=$False},
    @{Cell='D49'; Text='0.1233'; Numeric=$True},
    @{Cell='E49'; Text='  -2.91%  '; Numeric=$False},
    @{Cell='D50'; Text='34.89'; Numeric=$True},
    @{Cell='E50'; Text='  -0.23%  '; Numeric=$False},
    @{Cell='D51'; Text='0.8878'; Numeric=$True},
    @{Cell='E51'; Text='  +4.14%  '; Numeric=$False}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Numeric) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Text
}
